$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# --- Title ---
Replace-Text "Reforming Justice Systems for a Just Society" "Politics: Navigating the Maze of Power and Influence"

# --- Author name (3 runs collapse into 1) ---
Replace-Text "Isabella J. Hutchinson" "Emma Watson"

# --- Email address (keep the "." runs separate, only swap surrounding text) ---
Replace-Text "isabella" "emma"
Replace-Text "hutchinson@berkeley" "watson87@schoolmail"
Replace-Text "edu" "net"

# --- Body paragraph 1 (size 24) ---
Replace-Text "Justice systems worldwide strive to uphold fairness and ensure a safe and just society" "Politics, a multifaceted and dynamic realm of human interaction, permeates every aspect of our lives"
Replace-Text " However, disparities and inefficiencies persist, prompting calls for reforms to align justice systems with ideals of equality and impartiality" " It shapes the laws, policies, and decisions that govern societies, impacting individuals, communities, and nations alike"
Replace-Text " This essay delves into the pertinent issues confronting justice systems, encompassing the challenges and obstacles hindering fair outcomes, and explores viable solutions toward constructing a more just and equitable society for all" " As citizens of a democratic society, it is imperative for us to understand the intricacies of politics and the role we play in shaping its course"

Replace-Text "Justice systems are facing a multitude of challenges, ranging from racial and socioeconomic disparities to insufficient resources and outdated technologies" "Politics is often perceived as a complex web of power dynamics, negotiations, and compromises"
Replace-Text " These challenges manifest in unequal treatment, wrongful convictions, and disproportionate incarceration rates that disproportionately affect vulnerable and marginalized communities" " It involves the interactions among various stakeholders, including elected officials, political parties, interest groups, and the general public"
Replace-Text " The burden falls heavily on the shoulders of those navigating the justice system, eroding trust and fostering an atmosphere of injustice" " Understanding the different branches of government, their functions, and how they interact is crucial for comprehending the political landscape"

# Insert the two new sentences (new "." run + new sentence run) right after the sentence above,
# before the pre-existing "." run that used to close off "...atmosphere of injustice."
$rng = $d.Content
$rng.Find.Execute("Understanding the different branches of government, their functions, and how they interact is crucial for comprehending the political landscape") | Out-Null
$rng.Collapse(0)  # wdCollapseEnd
$rng.InsertAfter(".")
$rng.Collapse(0)
$rng.InsertAfter(" Political ideologies, such as liberalism, conservatism, and socialism, influence the policy positions and actions of political actors, and it is essential to grasp these ideologies and their implications")

Replace-Text "Furthermore, the increasing complexity of modern society presents new challenges that traditional justice systems may be ill-equipped to handle effectively" "Beyond the formal institutions and processes, politics also encompasses the informal dynamics of influence and persuasion"
Replace-Text " The rise of cybercrime, the prevalence of intellectual property disputes, and the vulnerabilities posed by globalization underscore the need for reforms that adapt justice systems to address evolving societal needs" " Lobbying, public relations, and grassroots movements play a significant role in shaping political outcomes"

# Remove the double line-break + "To rectify..." sentence, replacing with plain new sentence text (no breaks)
Replace-Text "To rectify these challenges and advance towards a more just society, comprehensive reforms are necessary" "placeholder-to-rectify"
$rng2 = $d.Content
$rng2.Find.Execute("placeholder-to-rectify") | Out-Null
# extend selection backward over the two breaks that precede this run
$rng2.MoveStart(1, -2) | Out-Null
$rng2.Text = " The media plays a vital role in informing and shaping public opinion, and understanding the relationship between politics and the media is crucial for informed citizenship"

Replace-Text " The implementation of restorative justice practices, which focus on healing and rehabilitation, could mitigate the harmful effects of mass incarceration while fostering accountability and seeking reparations for victims" " The influence of money in politics, campaign finance regulations, and the role of special interest groups are important aspects to consider in analyzing the political landscape"

# Remove the trailing ". Moreover, investing..." sentence entirely (was the final sentence in this paragraph)
$rng3 = $d.Content
$rng3.Find.Execute(" Moreover, investing in legal aid services and providing adequate resources for public defenders can help level the playing field, ensuring that all individuals have access to competent legal representation") | Out-Null
$rng3.MoveStart(1, -1) | Out-Null
$rng3.Text = ""

# --- Summary heading stays the same ---

# --- Final summary paragraph ---
Replace-Text "Justice systems are facing challenges that undermine the pursuit of fairness and equality" "Politics, a complex and ever-evolving field, involves the interactions among various stakeholders, encompassing formal institutions, informal dynamics, and the influence of the media"
Replace-Text " Racial and socioeconomic disparities, insufficient resources, and outdated technologies contribute to unequal treatment and wrongful convictions" " It requires an understanding of political ideologies, government structures, and the role of interest groups"
Replace-Text " As society evolves, new challenges emerge, such as cybercrime and intellectual property disputes, necessitating reforms to adapt justice systems" " As citizens, it is essential for us to engage with politics, be informed about current issues, and participate in the democratic process to create a society that reflects our values and aspirations"

# Remove the trailing three old sentences (". To promote justice...." + ". These measures....")
$rng4 = $d.Content
$rng4.Find.Execute(" To promote justice, reforms should incorporate restorative justice practices, expand access to legal aid, and leverage technology to improve efficiency and transparency") | Out-Null
$rng4.MoveEndUntil(".", 1000) | Out-Null
$rng4.MoveEnd(1, 1) | Out-Null
$rng4End = $d.Content
$rng4End.Find.Execute(" These measures, along with other comprehensive reforms, can pave the way for a just society where the scales of justice tip equitably for all") | Out-Null
$rng4.End = $rng4End.End
$rng4.Text = ""

# --- Add a trailing empty paragraph at the very end of the document body ---
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
